$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 251, pushing existing rows 251-286 down to 252-287.
$ws.Range("A251").EntireRow.Insert()

# Populate the newly inserted row 251 with the new observation's data.
$ws.Range("A251").Value = 6
$ws.Range("B251").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C251").Value = "Metropolitana"
$ws.Range("D251").Value = 44946
$ws.Range("D251").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E251").Value = 13
$ws.Range("F251").Value = 100112029
$ws.Range("G251").Value = "Orégano"
$ws.Range("H251").Value = "Sin especificar"
$ws.Range("I251").Value = "Primera"
$ws.Range("J251").Value = 48
$ws.Range("K251").Value = 19000
$ws.Range("L251").Value = 20000
$ws.Range("M251").Value = 19458
$ws.Range("N251").Value = '$/docena de atados'
$ws.Range("O251").Value = "Región Metropolitana"
$ws.Range("P251").Value = 6486
$ws.Range("Q251").Value = 3
$ws.Range("R251").Value = "Hortaliza"
